$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Examples" column content (row 2..16) ---
$examples = @(
  "If you procrastinate on this project, you will miss the deadline.",
  "He is such a procrastinator; he always studies the night before the exam.",
  "She is very organized; her desk is always clean and tidy.",
  "Remote work requires you to be self-motivated and disciplined.",
  "I know this is last minute, but can you help me review this report?",
  "Please send the files; I need it as soon as possible.",
  "Thanks for lending me your charger. You are a lifesaver.",
  "Fixing the printer? No sweat, I can do it in five minutes.",
  "I can help with the dishes; I know you've got a lot on your plate right now.",
  "I know you're busy, so I won't keep you any longer.",
  "Please sign these legal documents before leaving.",
  "Call me back asap (as soon as possible).",
  "We need a doctor here; it's really urgent.",
  "Thanks for covering my shift. I owe you one.",
  "The doctor will see you now. Sorry to keep you waiting."
)

for ($i = 0; $i -lt $examples.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 3).Value = $examples[$i]
}

# --- Fonts ---
# Header row (A1:C1) -> bold Arial 11 FF1F1F1F
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Name = "Arial"
$headerRange.Font.Size = 11
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 2039583

# Data rows (A2:C16) -> regular Arial 11 FF1F1F1F
$dataRange = $ws.Range("A2:C16")
$dataRange.Font.Name = "Arial"
$dataRange.Font.Size = 11
$dataRange.Font.Bold = $false
$dataRange.Font.Color = 2039583

# --- Borders: thin -> medium black on the whole used range ---
$allRange = $ws.Range("A1:C16")
$allRange.Borders.LineStyle = 1
$allRange.Borders.Weight = 4
$allRange.Borders.Color = 0

# --- Alignment: left / center / wrap / indent 1 / readingOrder ltr ---
$allRange.HorizontalAlignment = -4131  # xlLeft
$allRange.VerticalAlignment = -4108    # xlCenter
$allRange.WrapText = $true
$allRange.IndentLevel = 1
$allRange.ReadingOrder = 1             # xlLTR

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 29.88
$ws.Columns.Item(2).ColumnWidth = 35.31
$ws.Columns.Item(3).ColumnWidth = 55.88

# --- Row heights ---
$rowHeights = @{
  1 = 30.75
  2 = 29.25
  3 = 29.25
  4 = 15.75
  5 = 29.25
  6 = 29.25
  7 = 30.75
  8 = 15.75
  9 = 29.25
  10 = 30.75
  11 = 29.25
  12 = 15.75
  13 = 15.75
  14 = 15.75
  15 = 29.25
  16 = 15.75
}
foreach ($r in $rowHeights.Keys) {
  $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
}

# --- Selection ---
$ws.Range("A2:A16").Select()
